$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 86

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-28"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "19:24:37"
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "Wednesday"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "26"
$ws.Cells.Item($row, 5).Value = 123064
$ws.Cells.Item($row, 6).Value = 134275
$ws.Cells.Item($row, 7).Value = 163935
$ws.Cells.Item($row, 8).Value = 134135
$ws.Cells.Item($row, 9).Value = 177186
$ws.Cells.Item($row, 10).Value = 114972
$ws.Cells.Item($row, 11).Value = 204315
$ws.Cells.Item($row, 12).Value = 226470
$ws.Cells.Item($row, 13).Value = 176127
$ws.Cells.Item($row, 14).Value = 104406
$ws.Cells.Item($row, 15).Value = 39747
$ws.Cells.Item($row, 16).Value = 33720
$ws.Cells.Item($row, 17).Value = 52439
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36066
$ws.Cells.Item($row, 20).Value = -1
